$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.4320577737004783
$ws.Cells.Item(2, 4).Value = 0.3523579448139174
$ws.Cells.Item(2, 5).Value = 0.1821962538560413
$ws.Cells.Item(2, 6).Value = 5.009516939772851
$ws.Cells.Item(2, 7).Value = 0.002627872196580732
$ws.Cells.Item(2, 11).Value = 3.550942404564978
$ws.Cells.Item(2, 12).Value = 0.1587147059426144

$ws.Cells.Item(3, 3).Value = 0.4246816063167387
$ws.Cells.Item(3, 4).Value = 0.3440399599330135
$ws.Cells.Item(3, 5).Value = 0.1787022366166298
$ws.Cells.Item(3, 6).Value = 4.8349258661124
$ws.Cells.Item(3, 7).Value = 0.00263590399009006
$ws.Cells.Item(3, 11).Value = 3.42554775289409
$ws.Cells.Item(3, 12).Value = 0.155328413408192

$ws.Cells.Item(4, 3).Value = 0.4204215540843848
$ws.Cells.Item(4, 4).Value = 0.3389924247074987
$ws.Cells.Item(4, 5).Value = 0.1766738336087172
$ws.Cells.Item(4, 6).Value = 4.729635094922543
$ws.Cells.Item(4, 7).Value = 0.00264108197214553
$ws.Cells.Item(4, 11).Value = 3.351847938086394
$ws.Cells.Item(4, 12).Value = 0.153354004970268

$ws.Cells.Item(5, 3).Value = 0.4187526909531698
$ws.Cells.Item(5, 4).Value = 0.3369498019774255
$ws.Cells.Item(5, 5).Value = 0.1758764112900231
$ws.Cells.Item(5, 6).Value = 4.68719637647169
$ws.Cells.Item(5, 7).Value = 0.002643254278365769
$ws.Cells.Item(5, 11).Value = 3.322636396374605
$ws.Cells.Item(5, 12).Value = 0.1525755326534934

$ws.Cells.Item(6, 3).Value = 0.418479618802678
$ws.Cells.Item(6, 4).Value = 0.3366114677102132
$ws.Cells.Item(6, 5).Value = 0.1757457548130219
$ws.Cells.Item(6, 6).Value = 4.680177389013693
$ws.Cells.Item(6, 7).Value = 0.002643618754883724
$ws.Cells.Item(6, 11).Value = 3.317835291799497
$ws.Cells.Item(6, 12).Value = 0.1524478384464345

$ws.Cells.Item(7, 3).Value = 0.4203987759852907
$ws.Cells.Item(7, 4).Value = 0.3389648202922331
$ws.Cells.Item(7, 5).Value = 0.1766629614987174
$ws.Cells.Item(7, 6).Value = 4.729060871679934
$ws.Cells.Item(7, 7).Value = 0.002641111016278131
$ws.Cells.Item(7, 11).Value = 3.351450662049558
$ws.Cells.Item(7, 12).Value = 0.1533434007872074

$ws.Cells.Item(8, 3).Value = 0.4294583527084797
$ws.Cells.Item(8, 4).Value = 0.3494770588458351
$ws.Cells.Item(8, 5).Value = 0.1809671007699336
$ws.Cells.Item(8, 6).Value = 4.948914514041093
$ws.Cells.Item(8, 7).Value = 0.002630590576045824
$ws.Cells.Item(8, 11).Value = 3.507018948685243
$ws.Cells.Item(8, 12).Value = 0.1575252190310721

$ws.Cells.Item(9, 3).Value = 0.4493819527861547
$ws.Cells.Item(9, 4).Value = 0.3705997057779911
$ws.Cells.Item(9, 5).Value = 0.1903469065133336
$ws.Cells.Item(9, 6).Value = 5.395755375717499
$ws.Cells.Item(9, 7).Value = 0.002611902940949484
$ws.Cells.Item(9, 11).Value = 3.838532506724448
$ws.Cells.Item(9, 12).Value = 0.1665686232231423

$ws.Cells.Item(10, 3).Value = 0.4653711467032622
$ws.Cells.Item(10, 4).Value = 0.3864771988317273
$ws.Cells.Item(10, 5).Value = 0.197828229329069
$ws.Cells.Item(10, 6).Value = 5.734444667405768
$ws.Cells.Item(10, 7).Value = 0.00259934016109516
$ws.Cells.Item(10, 11).Value = 4.09869983276792
$ws.Cells.Item(10, 12).Value = 0.1737436726828747

$ws.Cells.Item(11, 3).Value = 0.4729466343567594
$ws.Cells.Item(11, 4).Value = 0.3937889490235875
$ws.Cells.Item(11, 5).Value = 0.2013637102903161
$ws.Cells.Item(11, 6).Value = 5.8909626584325
$ws.Cells.Item(11, 7).Value = 0.002593874674401277
$ws.Cells.Item(11, 11).Value = 4.220773146195711
$ws.Cells.Item(11, 12).Value = 0.177126915281832

$ws.Cells.Item(12, 3).Value = 0.4758593384651135
$ws.Cells.Item(12, 4).Value = 0.396571382986906
$ws.Cells.Item(12, 5).Value = 0.2027218258973917
$ws.Cells.Item(12, 6).Value = 5.950598255738896
$ws.Cells.Item(12, 7).Value = 0.002591840604884916
$ws.Cells.Item(12, 11).Value = 4.267543121386154
$ws.Cells.Item(12, 12).Value = 0.1784255203501601

$ws.Cells.Item(13, 3).Value = 0.4752300672189733
$ws.Cells.Item(13, 4).Value = 0.395971516288796
$ws.Cells.Item(13, 5).Value = 0.2024284679533963
$ws.Cells.Item(13, 6).Value = 5.937738178250697
$ws.Cells.Item(13, 7).Value = 0.002592277099571027
$ws.Cells.Item(13, 11).Value = 4.257446065126771
$ws.Cells.Item(13, 12).Value = 0.1781450619142504

$ws.Cells.Item(14, 3).Value = 0.4731853782311362
$ws.Cells.Item(14, 4).Value = 0.394017583485379
$ws.Cells.Item(14, 5).Value = 0.2014750546307056
$ws.Cells.Item(14, 6).Value = 5.895861506971642
$ws.Cells.Item(14, 7).Value = 0.002593706618569368
$ws.Cells.Item(14, 11).Value = 4.224610000073937
$ws.Cells.Item(14, 12).Value = 0.1772334009810237

$ws.Cells.Item(15, 3).Value = 0.4719386999514938
$ws.Cells.Item(15, 4).Value = 0.3928225436386583
$ws.Cells.Item(15, 5).Value = 0.2008935845792905
$ws.Cells.Item(15, 6).Value = 5.870258867517748
$ws.Cells.Item(15, 7).Value = 0.002594586867320042
$ws.Cells.Item(15, 11).Value = 4.204567986931238
$ws.Cells.Item(15, 12).Value = 0.1766772629603821

$ws.Cells.Item(16, 3).Value = 0.4648821983905975
$ws.Cells.Item(16, 4).Value = 0.3860012207119894
$ws.Cells.Item(16, 5).Value = 0.1975998625142452
$ws.Cells.Item(16, 6).Value = 5.724266254472099
$ws.Cells.Item(16, 7).Value = 0.002599702338602876
$ws.Cells.Item(16, 11).Value = 4.090797647649595
$ws.Cells.Item(16, 12).Value = 0.1735249948010136

$ws.Cells.Item(17, 3).Value = 0.4606310247598628
$ws.Cells.Item(17, 4).Value = 0.3818399446480782
$ws.Cells.Item(17, 5).Value = 0.1956133361986971
$ws.Cells.Item(17, 6).Value = 5.635340880198385
$ws.Cells.Item(17, 7).Value = 0.002602904196448297
$ws.Cells.Item(17, 11).Value = 4.02196227613922
$ws.Cells.Item(17, 12).Value = 0.1716219355906361

$ws.Cells.Item(18, 3).Value = 0.4582142426361884
$ws.Cells.Item(18, 4).Value = 0.3794548149926982
$ws.Cells.Item(18, 5).Value = 0.1944831628315313
$ws.Cells.Item(18, 6).Value = 5.584422793366798
$ws.Cells.Item(18, 7).Value = 0.002604769310060344
$ws.Cells.Item(18, 11).Value = 3.982719801239625
$ws.Cells.Item(18, 12).Value = 0.1705385551340868

$ws.Cells.Item(19, 3).Value = 0.4574008176782911
$ws.Cells.Item(19, 4).Value = 0.3786486559049536
$ws.Cells.Item(19, 5).Value = 0.1941026304826394
$ws.Cells.Item(19, 6).Value = 5.567221793232562
$ws.Cells.Item(19, 7).Value = 0.002605404847797444
$ws.Cells.Item(19, 11).Value = 3.969492806981464
$ws.Cells.Item(19, 12).Value = 0.1701736574805039

$ws.Cells.Item(20, 3).Value = 0.4610806277804045
$ws.Cells.Item(20, 4).Value = 0.3822820518454364
$ws.Cells.Item(20, 5).Value = 0.1958235173533538
$ws.Cells.Item(20, 6).Value = 5.644783293041058
$ws.Cells.Item(20, 7).Value = 0.002602560924158503
$ws.Cells.Item(20, 11).Value = 4.029253656072228
$ws.Cells.Item(20, 12).Value = 0.1718233573239871

$ws.Cells.Item(21, 3).Value = 0.4737847529569876
$ws.Cells.Item(21, 4).Value = 0.3945911242017246
$ws.Cells.Item(21, 5).Value = 0.2017545688730635
$ws.Cells.Item(21, 6).Value = 5.908151665783691
$ws.Cells.Item(21, 7).Value = 0.002593285770534317
$ws.Cells.Item(21, 11).Value = 4.234239931663922
$ws.Cells.Item(21, 12).Value = 0.1775007020823693

$ws.Cells.Item(22, 3).Value = 0.4823445807148801
$ws.Cells.Item(22, 4).Value = 0.4027156745735283
$ws.Cells.Item(22, 5).Value = 0.2057435264956808
$ws.Cells.Item(22, 6).Value = 6.082415912586725
$ws.Cells.Item(22, 7).Value = 0.002587431249989222
$ws.Cells.Item(22, 11).Value = 4.371382066011734
$ws.Cells.Item(22, 12).Value = 0.1813129970606582

$ws.Cells.Item(23, 3).Value = 0.4777523205895875
$ws.Cells.Item(23, 4).Value = 0.3983718625064512
$ws.Cells.Item(23, 5).Value = 0.203604134742676
$ws.Cells.Item(23, 6).Value = 5.989207557218265
$ws.Cells.Item(23, 7).Value = 0.0025905370363669
$ws.Cells.Item(23, 11).Value = 4.297893705178581
$ws.Cells.Item(23, 12).Value = 0.179268887831924

$ws.Cells.Item(24, 3).Value = 0.4608772773472651
$ws.Cells.Item(24, 4).Value = 0.3820821526828695
$ws.Cells.Item(24, 5).Value = 0.1957284573878582
$ws.Cells.Item(24, 6).Value = 5.640513736672915
$ws.Cells.Item(24, 7).Value = 0.002602716041815079
$ws.Cells.Item(24, 11).Value = 4.025956190039437
$ws.Cells.Item(24, 12).Value = 0.1717322612096268

$ws.Cells.Item(25, 3).Value = 0.4437574371574158
$ws.Cells.Item(25, 4).Value = 0.3648260540574171
$ws.Cells.Item(25, 5).Value = 0.1877070918920154
$ws.Cells.Item(25, 6).Value = 5.273107579000737
$ws.Cells.Item(25, 7).Value = 0.002616752237861843
$ws.Cells.Item(25, 11).Value = 3.745969289330958
$ws.Cells.Item(25, 12).Value = 0.1640301751885147

Write-Host "Updated pl_mw values for 380 kV case"